$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Objects_Alerts")
Write-Host $ws.Name
